$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 20, pushing existing row 20..101 down to 21..102
$ws.Rows(20).Insert()

# Populate the newly inserted row 20 with the new record
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C20").Value = "Ñuble"
$ws.Range("D20").Value = 44859
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = 100112021
$ws.Range("G20").Value = "Ají"
$ws.Range("H20").Value = "Inferno"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = 20000
$ws.Range("L20").Value = 20000
$ws.Range("M20").Value = 20000
$ws.Range("N20").Value = "$/caja 10 kilos"
$ws.Range("O20").Value = "Región de Arica y Parinacota"
$ws.Range("P20").Value = 2000
$ws.Range("Q20").Value = 10
$ws.Range("R20").Value = "Hortaliza"
